$d = $word.ActiveDocument
$d.Content.Find.Execute("2018./2019", $true, $false, $false, $false, $false, $true, 1, $false, "2019./2020", 2)
